$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update goal_error_height_not_met (row 26): new message + C value ---
$ws.Range('B26').Value = 'One or more builds do exceed the required height!'
$ws.Range('C26').Value = 2.5

# --- 2. Insert 8 rows at 37 for unit-cube related strings ---
$ws.Range('A37:A44').EntireRow.Insert()
$ws.Range('A37').Value = 'unit_cube'
$ws.Range('B37').Value = 'Unit Cube'
$ws.Range('A38').Value = 'length'
$ws.Range('B38').Value = 'Length'
$ws.Range('A39').Value = 'width'
$ws.Range('B39').Value = 'Width'
$ws.Range('A40').Value = 'height'
$ws.Range('B40').Value = 'Height'
$ws.Range('A41').Value = 'sides_eval'
$ws.Range('B41').Value = 'Sides: Length, Width, Height'
$ws.Range('A42').Value = 'volume_eval'
$ws.Range('B42').Value = 'Length x Width x Height = Volume'
$ws.Range('A43').Value = 'cubic_unit_eval'
$ws.Range('B43').Value = '1 Unit x 1 Unit x 1 Unit = 1 Cubic Unit'
$ws.Range('A44').Value = 'one_unit'
$ws.Range('B44').Value = '1 Unit'

# --- 3. Update next_instruct (now row 48) text ---
$ws.Range('B48').Value = 'Press this button to continue.'

# --- 4. Insert 5 rows at 49 for drag/verify/build instructs ---
$ws.Range('A49:A53').EntireRow.Insert()
$ws.Range('A49').Value = 'drag_material_instruct'
$ws.Range('B49').Value = 'Press and drag the material to the designated location.'
$ws.Range('A50').Value = 'drag_side_instruct'
$ws.Range('B50').Value = 'Press and drag the sides to expand the material.'
$ws.Range('A51').Value = 'expand_confirm_instruct'
$ws.Range('B51').Value = 'Press this button to finish expanding.'
$ws.Range('A52').Value = 'verify_instruct'
$ws.Range('B52').Value = 'Press here when you are ready to proceed.'
$ws.Range('A53').Value = 'build_instruct'
$ws.Range('B53').Value = 'Press this button to build.'

# --- 5. Update level_0_intro_0_3 (now row 56) text ---
$ws.Range('B56').Value = 'These buttons will allow you to rotate or elevate the view.'

# --- 6. Set level_0_intro_1_2 (now row 58) B value (previously empty) ---
$ws.Range('B58').Value = 'Remember that the volume of an object tells us how much it occupies space.'

# --- 7. Append 5 new rows (59-63) for more level_0_intro strings ---
$ws.Range('A59').Value = 'level_0_intro_1_3'
$ws.Range('B59').Value = 'By using unit cubes, we can easily determine where to place the materials, and how much volume we will need.'
$ws.Range('A60').Value = 'level_0_intro_1_4'
$ws.Range('B60').Value = 'So if an object is made up of unit cubes, then the volume is equal to the number of unit cubes that make up the object.'
$ws.Range('A61').Value = 'level_0_intro_1_5'
$ws.Range('B61').Value = 'A unit cube''s measurement can also be changed based on specific needs. In our case, one unit cube equals to one cubic feet.'
$ws.Range('A62').Value = 'level_0_intro_2_0'
$ws.Range('B62').Value = 'Here''s our first objective.'
$ws.Range('A63').Value = 'level_0_intro_2_1'
$ws.Range('B63').Value = 'You will be placing a number of unit cubes on the ground to match the required volume.'

# --- 8. Update view state: top-left cell and selection ---
$ws.Range('B50').Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
